$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.377.37'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.503.72'
$ws.Range('E3').Value = '  -3.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '200.45'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '551.90'
$ws.Range('E6').Value = '  -4.39%  '
$ws.Range('D7').Value = '3.498.23'
$ws.Range('E7').Value = '  -3.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').Value = '  -2.46%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.654'
$ws.Range('E10').Value = '  -3.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '63.45'
$ws.Range('E11').Value = '  +11.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.143'
$ws.Range('E12').Value = '  -7.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000269'
$ws.Range('E13').Value = '  -8.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.82'
$ws.Range('E14').Value = '  -3.32%  '
$ws.Range('D15').Value = '4.046.57'
$ws.Range('E15').Value = '  -4.11%  '
$ws.Range('D16').Value = '3.491.39'
$ws.Range('E16').Value = '  -4.06%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '67.045.59'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.30'
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.79'
$ws.Range('E20').Value = '  -6.31%  '
$ws.Range('E21').Value = '  -5.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.94'
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.17'
$ws.Range('E23').Value = '  -4.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.00'
$ws.Range('E24').Value = '  -5.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.34'
$ws.Range('E25').Value = '  -4.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.19'
$ws.Range('E27').Value = '  -3.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.80'
$ws.Range('E28').Value = '  -5.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.82'
$ws.Range('E29').Value = '  -3.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.96'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '677.63'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.99'
$ws.Range('E32').Value = '  -14.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.73'
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.93'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.110'
$ws.Range('E35').Value = '  -7.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.63'
$ws.Range('E36').Value = '  -9.80%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.397'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('D39').Value = '3.076.77'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('E41').Value = '  -4.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.97'
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('D43').Value = '0.0₃0673'
$ws.Range('E43').Value = '  -15.63%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.52'
$ws.Range('E44').Value = '  -12.16%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.77'
$ws.Range('E45').Value = '  +6.76%  '
$ws.Range('E46').Value = '  -9.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0397'
$ws.Range('E47').Value = '  -6.29%  '
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '136.87'
$ws.Range('E49').Value = '  -4.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.22'
$ws.Range('E50').Value = '  -7.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.88'
$ws.Range('E51').Value = '  -7.17%  '
